$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-10-29 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-10-30 Wednesday", 2)

$t = $d.Tables.Item(1)

# Each entry maps a specific table cell (1-based Row/Col) from its old
# text to its new text, listed in document order (top-to-bottom,
# left-to-right). Because earlier replacements can change text length and
# shift absolute character offsets, each cell's Range.Start/End is
# re-resolved live (via $t.Cell(...)) right before it is used, and cells
# are visited strictly in document order. This keeps the two identical
# "77÷6=12, 5" cells (row 13, cols 2 and 3) correctly distinguished even
# though their text content matches, and avoids any stale offsets from
# earlier edits.
$edits = @(
    @{Row=1;  Col=1; Old="19÷3=6, 1";   New="63÷6=10, 3"},
    @{Row=1;  Col=2; Old="15÷5=3, 0";   New="53÷3=17, 2"},
    @{Row=1;  Col=3; Old="70÷3=23, 1";  New="52÷8=6, 4"},
    @{Row=1;  Col=4; Old="99÷6=16, 3";  New="77÷9=8, 5"},
    @{Row=1;  Col=5; Old="49÷6=8, 1";   New="61÷8=7, 5"},

    @{Row=5;  Col=1; Old="88÷2=44, 0";  New="82÷9=9, 1"},
    @{Row=5;  Col=2; Old="12÷3=4, 0";   New="22÷9=2, 4"},
    @{Row=5;  Col=3; Old="36÷9=4, 0";   New="82÷4=20, 2"},
    @{Row=5;  Col=4; Old="36÷4=9, 0";   New="33÷3=11, 0"},
    @{Row=5;  Col=5; Old="91÷3=30, 1";  New="84÷7=12, 0"},

    @{Row=9;  Col=1; Old="17÷3=5, 2";   New="25÷3=8, 1"},
    @{Row=9;  Col=2; Old="89÷3=29, 2";  New="48÷6=8, 0"},
    @{Row=9;  Col=3; Old="79÷8=9, 7";   New="39÷9=4, 3"},
    @{Row=9;  Col=4; Old="49÷4=12, 1";  New="76÷3=25, 1"},
    @{Row=9;  Col=5; Old="18÷6=3, 0";   New="32÷6=5, 2"},

    @{Row=13; Col=1; Old="96÷3=32, 0";  New="52÷2=26, 0"},
    @{Row=13; Col=2; Old="77÷6=12, 5";  New="96÷2=48, 0"},
    @{Row=13; Col=3; Old="77÷6=12, 5";  New="78÷8=9, 6"},
    @{Row=13; Col=4; Old="60÷8=7, 4";   New="37÷2=18, 1"},
    @{Row=13; Col=5; Old="65÷3=21, 2";  New="19÷5=3, 4"},

    @{Row=17; Col=1; Old="13÷8=1, 5";   New="57÷6=9, 3"},
    @{Row=17; Col=2; Old="19÷2=9, 1";   New="80÷3=26, 2"},
    @{Row=17; Col=3; Old="24÷4=6, 0";   New="73÷5=14, 3"},
    @{Row=17; Col=4; Old="92÷2=46, 0";  New="90÷7=12, 6"},
    @{Row=17; Col=5; Old="29÷3=9, 2";   New="26÷8=3, 2"}
)

foreach ($e in $edits) {
    $cell = $t.Cell($e.Row, $e.Col)
    $r = $d.Range($cell.Range.Start, $cell.Range.End)
    $r.Find.Execute($e.Old, $true, $false, $false, $false, $false, $true, 1, $false, $e.New, 2) | Out-Null
}
